$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "id / label" rows (row 13, which is
# styled the same as the other bold-id / left-aligned-label rows 12-16)
# down onto the two new rows before filling in their values, so the new
# rows 18-19 keep the same look-and-feel as the rest of the table.
$ws.Range("A13:B13").Copy() | Out-Null
$ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:B13").Copy() | Out-Null
$ws.Range("A19:B19").PasteSpecial(-4122) | Out-Null

# New rows: id 17 -> "Misc. physical units", id 18 -> "Misc. units"
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Misc. physical units"
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Misc. units"

# Move the active selection down to A20, right below the newly added data,
# matching where the user's cursor ended up after the edit.
$ws.Range("A20").Select() | Out-Null
